# Auto-generated Excel COM-interop script
# Applies a batch of static-value updates (market price refresh) to the
# Kraken_Profits workbook, matching the scheduled-runner commit.

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 12
$ws.Range("H12").Value = 430
$ws.Range("J12").Value = 430
$ws.Range("L12").Value = 430
$ws.Range("N12").Value = -770
# Row 32
$ws.Range("H32").Value = 6908.5454
$ws.Range("I32").Value = 5250
$ws.Range("J32").Value = 7277.1113
$ws.Range("K32").Value = 5250
$ws.Range("L32").Value = 7277.1113
$ws.Range("M32").Value = -4924
$ws.Range("N32").Value = -7929.1113
# Row 61
$ws.Range("H61").Value = 3166
$ws.Range("I61").Value = 3166
$ws.Range("K61").Value = 9498
$ws.Range("M61").Value = -9326
# Row 70
$ws.Range("H70").Value = 25549.6
$ws.Range("J70").Value = 25549.6
$ws.Range("L70").Value = 76648.79999999999
$ws.Range("N70").Value = -77188.79999999999
# Row 73
$ws.Range("H73").Value = 25549.6
$ws.Range("J73").Value = 25549.6
$ws.Range("L73").Value = 76648.79999999999
$ws.Range("N73").Value = -78520.79999999999
# Row 86
$ws.Range("H86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("N86").Value = $null
# Row 89
$ws.Range("H89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("N89").Value = $null
# Row 112
$ws.Range("H112").Value = 7750
$ws.Range("J112").Value = 7750
$ws.Range("L112").Value = 23250
$ws.Range("N112").Value = -25466
# Row 113
$ws.Range("H113").Value = 1200
$ws.Range("I113").Value = 400
$ws.Range("K113").Value = 400
$ws.Range("M113").Value = 2854
# Row 137
$ws.Range("H137").Value = 2339.6
$ws.Range("J137").Value = 2232.6667
$ws.Range("L137").Value = 6698.000100000001
$ws.Range("N137").Value = -11798.0001
# Row 138
$ws.Range("H138").Value = 3482.6667
$ws.Range("J138").Value = 3999
$ws.Range("L138").Value = 11997
$ws.Range("N138").Value = -22277

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 50
$ws.Range("H50").Value = 1300
$ws.Range("I50").Value = 0
$ws.Range("J50").Value = 1300
$ws.Range("K50").Value = 0
$ws.Range("L50").Value = 1300
$ws.Range("M50").Value = $null
$ws.Range("N50").Value = -2728
# Row 61
$ws.Range("H61").Value = 3282.2
$ws.Range("I61").Value = 3282.2
$ws.Range("K61").Value = 3282.2
$ws.Range("M61").Value = -3070.2
# Row 63
$ws.Range("H63").Value = 926.4666999999999
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").Value = $null
# Row 66
$ws.Range("H66").Value = 926.4666999999999
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").Value = $null
# Row 122
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").Value = $null
# Row 136
$ws.Range("H136").Value = 3282.2
$ws.Range("I136").Value = 3282.2
$ws.Range("K136").Value = 9846.599999999999
$ws.Range("M136").Value = -7296.599999999999

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 23
$ws.Range("H23").Value = 496.5
$ws.Range("J23").Value = 496.5
$ws.Range("L23").Value = 496.5
$ws.Range("N23").Value = -1062.5
# Row 35
$ws.Range("H35").Value = 100
$ws.Range("I35").Value = 100
$ws.Range("K35").Value = 100
$ws.Range("M35").Value = 210

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 7
$ws.Range("H7").Value = 100.6
$ws.Range("I7").Value = 50.666668
$ws.Range("J7").Value = 175.5
$ws.Range("K7").Value = 50.666668
$ws.Range("L7").Value = 175.5
$ws.Range("M7").Value = 62.333332
$ws.Range("N7").Value = -401.5
# Row 22
$ws.Range("H22").Value = 728.1111
$ws.Range("I22").Value = 728.1111
$ws.Range("K22").Value = 728.1111
$ws.Range("M22").Value = -378.1111
# Row 50
$ws.Range("H50").Value = 27800
$ws.Range("J50").Value = 29500
$ws.Range("L50").Value = 29500
$ws.Range("N50").Value = -30750
# Row 82
$ws.Range("H82").Value = 59500
$ws.Range("I82").Value = 59000
$ws.Range("J82").Value = 60000
$ws.Range("K82").Value = 59000
$ws.Range("L82").Value = 60000
$ws.Range("M82").Value = -58639
$ws.Range("N82").Value = -60722
# Row 85
$ws.Range("H85").Value = 59500
$ws.Range("I85").Value = 59000
$ws.Range("J85").Value = 60000
$ws.Range("K85").Value = 59000
$ws.Range("L85").Value = 60000
$ws.Range("M85").Value = -57752
$ws.Range("N85").Value = -62496
# Row 122
$ws.Range("H122").Value = 471
$ws.Range("I122").Value = 471
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 1413
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = 1037
$ws.Range("N122").Value = $null

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 22
$ws.Range("H22").Value = 900
$ws.Range("I22").Value = 800
$ws.Range("J22").Value = 1000
$ws.Range("K22").Value = 2400
$ws.Range("L22").Value = 3000
$ws.Range("M22").Value = -2231
$ws.Range("N22").Value = -3338
# Row 27
$ws.Range("H27").Value = 900
$ws.Range("I27").Value = 800
$ws.Range("J27").Value = 1000
$ws.Range("K27").Value = 2400
$ws.Range("L27").Value = 3000
$ws.Range("M27").Value = -2298
$ws.Range("N27").Value = -3204
# Row 34
$ws.Range("H34").Value = 1123.5454
$ws.Range("I34").Value = 690.75
$ws.Range("J34").Value = 1370.8572
$ws.Range("K34").Value = 2072.25
$ws.Range("L34").Value = 4112.571599999999
$ws.Range("M34").Value = -1988.25
$ws.Range("N34").Value = -4280.571599999999
# Row 58
$ws.Range("H58").Value = 2250
$ws.Range("I58").Value = 2250
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 6750
$ws.Range("L58").Value = 0
$ws.Range("M58").Value = -6622
$ws.Range("N58").Value = $null
# Row 103
$ws.Range("H103").Value = 2227.4443
$ws.Range("I103").Value = 1198
$ws.Range("K103").Value = 3594
$ws.Range("M103").Value = -2715
# Row 137
$ws.Range("H137").Value = 500
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 500
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 1500
$ws.Range("M137").Value = $null
$ws.Range("N137").Value = -11700

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 13
$ws.Range("H13").Value = 1000
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 1000
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 1000
$ws.Range("M13").Value = $null
$ws.Range("N13").Value = -1278
# Row 23
$ws.Range("H23").Value = 4250
$ws.Range("J23").Value = 4250
$ws.Range("L23").Value = 4250
$ws.Range("N23").Value = -4696
# Row 102
$ws.Range("H102").Value = 2580.6924
$ws.Range("I102").Value = 2580.6924
$ws.Range("K102").Value = 2580.6924
$ws.Range("M102").Value = -958.6923999999999
# Row 122
$ws.Range("H122").Value = 1879
$ws.Range("I122").Value = 1750
$ws.Range("J122").Value = 2008
$ws.Range("K122").Value = 5250
$ws.Range("L122").Value = 6024
$ws.Range("M122").Value = -2800
$ws.Range("N122").Value = -10924
# Row 136
$ws.Range("H136").Value = 40000
$ws.Range("J136").Value = 40000
$ws.Range("L136").Value = 120000
$ws.Range("N136").Value = -125100

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 776.88464
$ws.Range("I22").Value = 777.55
$ws.Range("J22").Value = 774.6667
$ws.Range("K22").Value = 777.55
$ws.Range("L22").Value = 774.6667
$ws.Range("M22").Value = -482.55
$ws.Range("N22").Value = -1364.6667
# Row 27
$ws.Range("H27").Value = 776.88464
$ws.Range("I27").Value = 777.55
$ws.Range("J27").Value = 774.6667
$ws.Range("K27").Value = 777.55
$ws.Range("L27").Value = 774.6667
$ws.Range("M27").Value = -670.55
$ws.Range("N27").Value = -988.6667
# Row 40
$ws.Range("H40").Value = 5090.2856
$ws.Range("I40").Value = 5090.2856
$ws.Range("K40").Value = 5090.2856
$ws.Range("M40").Value = -4954.2856
# Row 55
$ws.Range("H55").Value = 3022.2307
$ws.Range("I55").Value = 2216.5
$ws.Range("J55").Value = 3712.8572
$ws.Range("K55").Value = 2216.5
$ws.Range("L55").Value = 3712.8572
$ws.Range("M55").Value = -2043.5
$ws.Range("N55").Value = -4058.8572
# Row 82
$ws.Range("H82").Value = 1961.4
$ws.Range("I82").Value = 1818.8572
$ws.Range("J82").Value = 2086.125
$ws.Range("K82").Value = 1818.8572
$ws.Range("L82").Value = 2086.125
$ws.Range("M82").Value = -1457.8572
$ws.Range("N82").Value = -2808.125
# Row 85
$ws.Range("H85").Value = 1961.4
$ws.Range("I85").Value = 1818.8572
$ws.Range("J85").Value = 2086.125
$ws.Range("K85").Value = 1818.8572
$ws.Range("L85").Value = 2086.125
$ws.Range("M85").Value = -570.8571999999999
$ws.Range("N85").Value = -4582.125
# Row 93
$ws.Range("H93").Value = 500
$ws.Range("I93").Value = 500
$ws.Range("K93").Value = 500
$ws.Range("M93").Value = 748
# Row 101
$ws.Range("H101").Value = 15751
$ws.Range("J101").Value = 15751
$ws.Range("L101").Value = 15751
$ws.Range("N101").Value = -22241
# Row 122
$ws.Range("H122").Value = 7000
$ws.Range("I122").Value = 9500
$ws.Range("J122").Value = 4500
$ws.Range("K122").Value = 28500
$ws.Range("L122").Value = 13500
$ws.Range("M122").Value = -26050
$ws.Range("N122").Value = -18400

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 15
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 0
$ws.Range("K15").Value = 0
$ws.Range("M15").Value = $null
# Row 82
$ws.Range("H82").Value = 0
$ws.Range("J82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("N82").Value = $null
# Row 85
$ws.Range("H85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("N85").Value = $null
# Row 113
$ws.Range("H113").Value = 727.5454999999999
$ws.Range("I113").Value = 723
$ws.Range("J113").Value = 735.5
$ws.Range("K113").Value = 2169
$ws.Range("L113").Value = 2206.5
$ws.Range("M113").Value = 1
$ws.Range("N113").Value = -6546.5
# Row 122
$ws.Range("H122").Value = 1983.625
$ws.Range("I122").Value = 2314.1667
$ws.Range("K122").Value = 6942.500100000001
$ws.Range("M122").Value = -4492.500100000001
# Row 126
$ws.Range("H126").Value = 3005.5
$ws.Range("I126").Value = 2551.4546
$ws.Range("K126").Value = 7654.3638
$ws.Range("M126").Value = -5184.3638
